$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 73: IVAN / UAPRENDIZAJE - se quitan campos semestre_sug y caracter de los triggers
$ws.Cells.Item(73, 1).Value = "IVAN"
$ws.Cells.Item(73, 2).Value = "UAPRENDIZAJE"
$ws.Cells.Item(73, 3).Value = "Se modificaron lo triggers de la bitacora quitando los campos semestre_sug y carácter"
$ws.Cells.Item(73, 4).Value = 41948

# Fila 74: IVAN / P_UA - se agregan campos semestre_sug y caracter a los triggers
$ws.Cells.Item(74, 1).Value = "IVAN"
$ws.Cells.Item(74, 2).Value = "P_UA"
$ws.Cells.Item(74, 3).Value = "Se modificaron lo triggers de la bitacora añadiendo los campos semestre_sug y carácter"
$ws.Cells.Item(74, 4).Value = 41949

# Replicar formato de las filas anteriores (centrado horizontal/vertical en A:B, fecha centrada en D)
$dateFormat = $ws.Range("D69").NumberFormat
foreach ($r in 73..74) {
    $ws.Range("A" + $r + ":B" + $r).HorizontalAlignment = -4108
    $ws.Range("A" + $r + ":B" + $r).VerticalAlignment = -4108
    $ws.Range("D" + $r).NumberFormat = $dateFormat
    $ws.Range("D" + $r).HorizontalAlignment = -4108
}

# Ajustar la vista: desplazar el panel dividido (topLeftCell) y la selección activa
$excel.ActiveWindow.SplitRow = 63
$ws.Range("A75").Select()
